# Remove the two empty "Unnamed" columns (D and E). Deleting column D twice
# shifts the old column F ("Dia da semana") into column D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(4).Delete()
$ws.Columns.Item(4).Delete()

# Append the three new data rows (382-384) that were added at the bottom of
# the sheet.
function Set-TextCell($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell "A382" "03/28/2021"
$ws.Range("B382").Value = 0.51
$ws.Range("C382").Value = 0.51
Set-TextCell "D382" "Domingo"

Set-TextCell "A383" "03/29/2021"
$ws.Range("B383").Value = 0.44
$ws.Range("C383").Value = 0.42
Set-TextCell "D383" "Segunda-feira"

Set-TextCell "A384" "03/30/2021"
$ws.Range("B384").Value = 0.44
$ws.Range("C384").Value = 0.42
Set-TextCell "D384" "Terça-feira"
